# Swap the contents of columns E and F (header labels "soft"/"rigid" and
# all the numeric counts below them) across the whole used range, per the
# "new layout G3 graphics" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row  # xlUp
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)  # column E
    $fCell = $ws.Cells.Item($r, 6)  # column F

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null -and $fVal -eq $null) {
        continue
    }

    $eCell.Value2 = $fVal
    $fCell.Value2 = $eVal
}
